# Rectify a year error
# - Insert a new worksheet "Mbt2018-mean-ste" right before "Mbt2018-mean"
#   with its own small results table (including the new FID/LPIPS/DISTS
#   metric columns).
# - Rename "FLIC2025" to "FTIC2024" (the year in that sheet's name was wrong).

$wb = $excel.ActiveWorkbook

# --- Insert the new sheet in the correct tab position -----------------
$anchor = $wb.Worksheets.Item("Mbt2018-mean")
$newSheet = $wb.Worksheets.Add($anchor)
$newSheet.Name = "Mbt2018-mean-ste"

# --- Header row ---------------------------------------------------------
$newSheet.Range("A1").Value = "λ"
$newSheet.Range("B1").Value = "bpp"
$newSheet.Range("C1").Value = "bpp-est"
$newSheet.Range("D1").Value = "Δbpp"
$newSheet.Range("E1").Value = "psnr"
$newSheet.Range("F1").Value = "psnr-est"
$newSheet.Range("G1").Value = "Δpsnr"
$newSheet.Range("H1").Value = "ms-ssim"
$newSheet.Range("I1").Value = "ms-ssim-est"
$newSheet.Range("J1").Value = "Δms-ssim"
$newSheet.Range("K1").Value = "FID"
$newSheet.Range("L1").Value = "LPIPS"
$newSheet.Range("M1").Value = "DISTS"

# --- Data rows ------------------------------------------------------------
$newSheet.Range("A2").Value = 0.0016000000000000001
$newSheet.Range("B2").Value = 0.159104
$newSheet.Range("C2").Value = 0.30343994311988348
$newSheet.Range("D2").Value = 0.14433594311988349
$newSheet.Range("E2").Value = 28.143000046412151
$newSheet.Range("F2").Value = 28.32690278689067
$newSheet.Range("G2").Value = 0.1839027404785156
$newSheet.Range("H2").Value = 0.92384535074234009
$newSheet.Range("I2").Value = 0.92540520429611206
$newSheet.Range("J2").Value = 0.0015598535537719731
$newSheet.Range("K2").Value = 110.49837618713281
$newSheet.Range("L2").Value = 0.32268826166788739
$newSheet.Range("M2").Value = 0.37496660898129153

$newSheet.Range("A3").Value = 0.0016100000000000001
$newSheet.Range("B3").Value = 0.15359999999999999
$newSheet.Range("C3").Value = 0.29559288173913961
$newSheet.Range("D3").Value = 0.14199288173913949
$newSheet.Range("E3").Value = 28.229444265365601
$newSheet.Range("F3").Value = 27.979055484135941
$newSheet.Range("G3").Value = 0.25038878122965608
$newSheet.Range("H3").Value = 0.92492032051086426
$newSheet.Range("I3").Value = 0.92181806017955148
$newSheet.Range("J3").Value = 0.0031022603313127779
$newSheet.Range("K3").Value = 109.6410637837439
$newSheet.Range("L3").Value = 0.31777926969031489
$newSheet.Range("M3").Value = 0.37836939593156182

# --- Fix the year typo in the last tab's name --------------------------
$wb.Worksheets.Item("FLIC2025").Name = "FTIC2024"

# --- Restore cursor positions / active tab -----------------------------
$newSheet.Range("E17").Select()
$wb.Worksheets.Item("WeConvene2024").Range("P15").Select()
$wb.Worksheets.Item("FTIC2024").Range("W47").Select()
